$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-11 21:48:43"
$ws.Range("O2").Value = "2.8 °C"
$ws.Range("E3").Value = "2026-02-11 21:48:45"
$ws.Range("H3").Value = "'80%"
$ws.Range("L3").Value = "94.7 km/h - 227º 21:05 TU"
$ws.Range("O3").Value = "0.5 °C"
$ws.Range("E4").Value = "2026-02-11 21:48:48"
$ws.Range("J4").Value = "1001.9 hPa"
$ws.Range("O4").Value = "15.9 °C"
$ws.Range("E5").Value = "2026-02-11 21:48:51"
$ws.Range("E6").Value = "2026-02-11 21:48:54"
$ws.Range("J6").Value = "1002.3 hPa"
$ws.Range("E7").Value = "2026-02-11 21:48:56"
$ws.Range("J7").Value = "1003.2 hPa"
$ws.Range("N7").Value = "15.0 °C 21:22 TU"
$ws.Range("O7").Value = "18.8 °C"
$ws.Range("E8").Value = "2026-02-11 21:48:59"
$ws.Range("H8").Value = "'59%"
$ws.Range("O8").Value = "14.8 °C"
$ws.Range("E9").Value = "2026-02-11 21:49:01"
$ws.Range("H9").Value = "'89%"
$ws.Range("I9").Value = "3.4 mm"
$ws.Range("E10").Value = "2026-02-11 21:49:04"
$ws.Range("H10").Value = "'76%"
$ws.Range("E11").Value = "2026-02-11 21:49:06"
$ws.Range("H11").Value = "'82%"
$ws.Range("I11").Value = "1.7 mm"
$ws.Range("E12").Value = "2026-02-11 21:49:09"
$ws.Range("I12").Value = "10.5 mm"
$ws.Range("E13").Value = "2026-02-11 21:49:11"
$ws.Range("I13").Value = "2.3 mm"
$ws.Range("J13").Value = "1004.9 hPa"
$ws.Range("E14").Value = "2026-02-11 21:49:14"
$ws.Range("H14").Value = "'52%"
$ws.Range("N14").Value = "13.0 °C 21:29 TU"
$ws.Range("O14").Value = "18.5 °C"
$ws.Range("E15").Value = "2026-02-11 21:49:17"
$ws.Range("I15").Value = "3.7 mm"
$ws.Range("E16").Value = "2026-02-11 21:49:19"
$ws.Range("H16").Value = "'68%"
$ws.Range("I16").Value = "10.6 mm"
$ws.Range("E17").Value = "2026-02-11 21:49:22"
$ws.Range("E18").Value = "2026-02-11 21:49:25"
$ws.Range("J18").Value = "1002.5 hPa"
$ws.Range("E19").Value = "2026-02-11 21:49:27"
$ws.Range("I19").Value = "2.7 mm"
$ws.Range("O19").Value = "9.0 °C"
$ws.Range("E20").Value = "2026-02-11 21:49:30"
$ws.Range("I20").Value = "4.5 mm"
$ws.Range("E21").Value = "2026-02-11 21:49:33"
$ws.Range("I21").Value = "4.7 mm"
$ws.Range("J21").Value = "1005.3 hPa"
$ws.Range("O21").Value = "8.2 °C"
$ws.Range("E22").Value = "2026-02-11 21:49:35"
$ws.Range("L22").Value = "113.0 km/h - 295º 21:11 TU"
$ws.Range("O22").Value = "-2.4 °C"
$ws.Range("E23").Value = "2026-02-11 21:49:38"
$ws.Range("E24").Value = "2026-02-11 21:49:40"
$ws.Range("J24").Value = "1006.4 hPa"
$ws.Range("O24").Value = "13.0 °C"
$ws.Range("E25").Value = "2026-02-11 21:49:43"
$ws.Range("H25").Value = "'67%"
$ws.Range("N25").Value = "-0.5 °C 21:04 TU"
$ws.Range("O25").Value = "1.6 °C"
$ws.Range("E26").Value = "2026-02-11 21:49:46"
$ws.Range("H26").Value = "'68%"
$ws.Range("J26").Value = "1002.3 hPa"
$ws.Range("E27").Value = "2026-02-11 21:49:49"
$ws.Range("E28").Value = "2026-02-11 21:49:51"
$ws.Range("J28").Value = "1002.6 hPa"
$ws.Range("E29").Value = "2026-02-11 21:49:54"
$ws.Range("E30").Value = "2026-02-11 21:49:57"
$ws.Range("I30").Value = "6.5 mm"
$ws.Range("J30").Value = "1002.5 hPa"
$ws.Range("O30").Value = "11.8 °C"
$ws.Range("E31").Value = "2026-02-11 21:50:00"
$ws.Range("H31").Value = "'66%"
$ws.Range("I31").Value = "3.5 mm"
$ws.Range("J31").Value = "1001.7 hPa"
$ws.Range("E32").Value = "2026-02-11 21:50:02"
$ws.Range("E33").Value = "2026-02-11 21:50:05"
$ws.Range("J33").Value = "1004.5 hPa"
$ws.Range("E34").Value = "2026-02-11 21:50:08"
$ws.Range("E35").Value = "2026-02-11 21:50:11"
$ws.Range("E36").Value = "2026-02-11 21:50:13"
$ws.Range("I36").Value = "8.8 mm"
$ws.Range("J36").Value = "1002.7 hPa"
$ws.Range("K36").Value = "10.5 MJ/m2"
$ws.Range("E37").Value = "2026-02-11 21:50:16"
$ws.Range("I37").Value = "1.5 mm"
$ws.Range("O37").Value = "9.3 °C"
$ws.Range("E38").Value = "2026-02-11 21:50:19"
$ws.Range("H38").Value = "'61%"
$ws.Range("O38").Value = "15.6 °C"
$ws.Range("E39").Value = "2026-02-11 21:50:21"
$ws.Range("H39").Value = "'59%"
$ws.Range("I39").Value = "4.0 mm"
$ws.Range("E40").Value = "2026-02-11 21:50:24"
$ws.Range("I40").Value = "7.6 mm"
$ws.Range("J40").Value = "1006.5 hPa"
$ws.Range("E41").Value = "2026-02-11 21:50:27"
$ws.Range("H41").Value = "'52%"
$ws.Range("J41").Value = "1004.3 hPa"
$ws.Range("K41").Value = "8.9 MJ/m2"
$ws.Range("O41").Value = "18.4 °C"
$ws.Range("E42").Value = "2026-02-11 21:50:30"
$ws.Range("O42").Value = "13.0 °C"
$ws.Range("E43").Value = "2026-02-11 21:50:32"
$ws.Range("H43").Value = "'67%"
$ws.Range("I43").Value = "8.2 mm"
$ws.Range("N43").Value = "9.5 °C 21:25 TU"
$ws.Range("O43").Value = "12.8 °C"
$ws.Range("E44").Value = "2026-02-11 21:50:35"
$ws.Range("L44").Value = "111.6 km/h - 202º 21:07 TU"
$ws.Range("E45").Value = "2026-02-11 21:50:38"
$ws.Range("J45").Value = "1005.1 hPa"
$ws.Range("E46").Value = "2026-02-11 21:50:41"
$ws.Range("I46").Value = "4.6 mm"
$ws.Range("J46").Value = "1006.7 hPa"
$ws.Range("N46").Value = "12.2 °C 21:00 TU"
$ws.Range("O46").Value = "16.6 °C"
